$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3555.5
$ws.Range("J76").Value = 3666.6667
$ws.Range("L76").Value = 3666.6667
$ws.Range("N76").Value = -4296.6667

$ws.Range("H79").Value = 3555.5
$ws.Range("J79").Value = 3666.6667
$ws.Range("L79").Value = 3666.6667
$ws.Range("N79").Value = -5850.6667

$ws.Range("H86").Value = 6909.55
$ws.Range("I86").Value = 2016.5385
$ws.Range("J86").Value = 15996.571
$ws.Range("K86").Value = 2016.5385
$ws.Range("L86").Value = 15996.571
$ws.Range("M86").Value = -893.5385000000001
$ws.Range("N86").Value = -18242.571

$ws.Range("H89").Value = 6909.55
$ws.Range("I89").Value = 2016.5385
$ws.Range("J89").Value = 15996.571
$ws.Range("K89").Value = 10082.6925
$ws.Range("L89").Value = 79982.855
$ws.Range("M89").Value = -4466.692500000001
$ws.Range("N89").Value = -91214.855

$ws.Range("H92").Value = 111111970
$ws.Range("I92").Value = 125000950
$ws.Range("J92").Value = 98
$ws.Range("K92").Value = 125000950
$ws.Range("L92").Value = 98
$ws.Range("M92").Value = -124999702
$ws.Range("N92").Value = -2594

$ws.Range("H96").Value = 16666962
$ws.Range("J96").Value = 180.33333
$ws.Range("L96").Value = 540.99999
$ws.Range("N96").Value = -3286.99999

$ws.Range("H125").Value = 1065.8889
$ws.Range("J125").Value = 1136.625
$ws.Range("L125").Value = 10229.625
$ws.Range("N125").Value = -15149.625

$ws.Range("H129").Value = 1160.2041
$ws.Range("J129").Value = 1273.8096
$ws.Range("L129").Value = 3821.4288
$ws.Range("N129").Value = -13821.4288

$ws.Range("H132").Value = 1747.305
$ws.Range("I132").Value = 1832.5636
$ws.Range("J132").Value = 575
$ws.Range("K132").Value = 5497.6908
$ws.Range("L132").Value = 1725
$ws.Range("M132").Value = -2967.6908
$ws.Range("N132").Value = -6785

$ws.Range("H137").Value = 1293.9166
$ws.Range("I137").Value = 1183.907
$ws.Range("K137").Value = 3551.721
$ws.Range("M137").Value = -1001.721

$ws.Range("H138").Value = 2610
$ws.Range("J138").Value = 2395.577
$ws.Range("L138").Value = 7186.731000000001
$ws.Range("N138").Value = -17466.731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 16

$ws.Range("H5").Value = 243
$ws.Range("I5").Value = 259.8
$ws.Range("J5").Value = 201
$ws.Range("K5").Value = 259.8
$ws.Range("L5").Value = 201
$ws.Range("M5").Value = -147.8
$ws.Range("N5").Value = -425

$ws.Range("H32").Value = 4722.4683
$ws.Range("I32").Value = 5522.237
$ws.Range("J32").Value = 1345.6666
$ws.Range("K32").Value = 5522.237
$ws.Range("L32").Value = 1345.6666
$ws.Range("M32").Value = -5235.237
$ws.Range("N32").Value = -1919.6666

$ws.Range("H45").Value = 4309.3335
$ws.Range("I45").Value = 4932
$ws.Range("J45").Value = 3913.0908
$ws.Range("K45").Value = 4932
$ws.Range("L45").Value = 3913.0908
$ws.Range("M45").Value = -4555
$ws.Range("N45").Value = -4667.0908

$ws.Range("H63").Value = 3544.111
$ws.Range("I63").Value = 2799.8
$ws.Range("J63").Value = 4474.5
$ws.Range("K63").Value = 2799.8
$ws.Range("L63").Value = 4474.5
$ws.Range("M63").Value = -2113.8
$ws.Range("N63").Value = -5846.5

$ws.Range("H66").Value = 3544.111
$ws.Range("I66").Value = 2799.8
$ws.Range("J66").Value = 4474.5
$ws.Range("K66").Value = 13999
$ws.Range("L66").Value = 22372.5
$ws.Range("M66").Value = -10567
$ws.Range("N66").Value = -29236.5

$ws.Range("H97").Value = 1413.8422
$ws.Range("I97").Value = 1403.25
$ws.Range("J97").Value = 1470.3334
$ws.Range("K97").Value = 1403.25
$ws.Range("L97").Value = 1470.3334
$ws.Range("M97").Value = -907.25
$ws.Range("N97").Value = -2462.3334

$ws.Range("H101").Value = 50049.75
$ws.Range("J101").Value = 50049.75
$ws.Range("L101").Value = 50049.75
$ws.Range("N101").Value = -56539.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 243
$ws.Range("I4").Value = 259.8
$ws.Range("J4").Value = 201
$ws.Range("K4").Value = 259.8
$ws.Range("L4").Value = 201
$ws.Range("M4").Value = -144.8
$ws.Range("N4").Value = -431

$ws.Range("H86").Value = 1767.8966
$ws.Range("I86").Value = 1531.5714
$ws.Range("J86").Value = 2388.25
$ws.Range("K86").Value = 1531.5714
$ws.Range("L86").Value = 2388.25
$ws.Range("M86").Value = -408.5714
$ws.Range("N86").Value = -4634.25

$ws.Range("H89").Value = 1767.8966
$ws.Range("I89").Value = 1531.5714
$ws.Range("J89").Value = 2388.25
$ws.Range("K89").Value = 7657.857
$ws.Range("L89").Value = 11941.25
$ws.Range("M89").Value = -2041.857
$ws.Range("N89").Value = -23173.25

$ws.Range("H134").Value = 3001
$ws.Range("I134").Value = 3256.7896
$ws.Range("K134").Value = 9770.3688
$ws.Range("M134").Value = -7235.3688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2854.889
$ws.Range("I31").Value = 1948.2727
$ws.Range("K31").Value = 1948.2727
$ws.Range("M31").Value = -1653.2727

$ws.Range("H34").Value = 2854.889
$ws.Range("I34").Value = 1948.2727
$ws.Range("K34").Value = 1948.2727
$ws.Range("M34").Value = -1746.2727

$ws.Range("H68").Value = 40206.332
$ws.Range("J68").Value = 40206.332
$ws.Range("L68").Value = 40206.332
$ws.Range("N68").Value = -41704.332

$ws.Range("H71").Value = 40206.332
$ws.Range("J71").Value = 40206.332
$ws.Range("L71").Value = 120618.996
$ws.Range("N71").Value = -128106.996

$ws.Range("H132").Value = 2337.0417
$ws.Range("I132").Value = 1758.1714
$ws.Range("K132").Value = 5274.5142
$ws.Range("M132").Value = -2744.5142

$ws.Range("H134").Value = 907.8889
$ws.Range("I134").Value = 786.8095
$ws.Range("K134").Value = 2360.4285
$ws.Range("M134").Value = 174.5715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 816.6667
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 840
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 2520
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -4142

$ws.Range("H71").Value = 816.6667
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 840
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 7560
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -15672

$ws.Range("H107").Value = 4582.7393
$ws.Range("I107").Value = 16833.166
$ws.Range("J107").Value = 259.05884
$ws.Range("K107").Value = 50499.49800000001
$ws.Range("L107").Value = 777.17652
$ws.Range("M107").Value = -48579.49800000001
$ws.Range("N107").Value = -4617.17652

$ws.Range("H113").Value = 613.7222
$ws.Range("I113").Value = 598.0909
$ws.Range("J113").Value = 638.2857
$ws.Range("K113").Value = 1794.2727
$ws.Range("L113").Value = 1914.8571
$ws.Range("M113").Value = 375.7273
$ws.Range("N113").Value = -6254.8571

$ws.Range("H131").Value = 794.5106
$ws.Range("J131").Value = 800.5056
$ws.Range("L131").Value = 2401.5168
$ws.Range("N131").Value = -12481.5168

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4993.8276
$ws.Range("I126").Value = 3806.3684
$ws.Range("K126").Value = 11419.1052
$ws.Range("M126").Value = -8949.1052

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4244.5557
$ws.Range("J22").Value = 5250
$ws.Range("L22").Value = 5250
$ws.Range("N22").Value = -5840

$ws.Range("H27").Value = 4244.5557
$ws.Range("J27").Value = 5250
$ws.Range("L27").Value = 5250
$ws.Range("N27").Value = -5464

$ws.Range("H40").Value = 6232
$ws.Range("J40").Value = 6801.091
$ws.Range("L40").Value = 6801.091
$ws.Range("N40").Value = -7073.091

$ws.Range("H93").Value = 3418.9092
$ws.Range("I93").Value = 3614.2856
$ws.Range("J93").Value = 3077
$ws.Range("K93").Value = 3614.2856
$ws.Range("L93").Value = 3077
$ws.Range("M93").Value = -2366.2856
$ws.Range("N93").Value = -5573

$ws.Range("H97").Value = 20000
$ws.Range("J97").Value = 20000
$ws.Range("L97").Value = 20000
$ws.Range("N97").Value = -21982

$ws.Range("H136").Value = 1333.2667
$ws.Range("I136").Value = 1199.9166
$ws.Range("J136").Value = 1866.6666
$ws.Range("K136").Value = 3599.7498
$ws.Range("L136").Value = 5599.9998
$ws.Range("M136").Value = -1049.7498
$ws.Range("N136").Value = -10699.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2739.818
$ws.Range("I62").Value = 2142.5
$ws.Range("J62").Value = 4332.6665
$ws.Range("K62").Value = 2142.5
$ws.Range("L62").Value = 4332.6665
$ws.Range("M62").Value = -1518.5
$ws.Range("N62").Value = -5580.6665

$ws.Range("H65").Value = 2739.818
$ws.Range("I65").Value = 2142.5
$ws.Range("J65").Value = 4332.6665
$ws.Range("K65").Value = 10712.5
$ws.Range("L65").Value = 21663.3325
$ws.Range("M65").Value = -7592.5
$ws.Range("N65").Value = -27903.3325

$ws.Range("H96").Value = 3000
$ws.Range("I96").Value = 3000
$ws.Range("K96").Value = 3000
$ws.Range("M96").Value = -1627

$ws.Range("H136").Value = 29413126
$ws.Range("I136").Value = 35715532
$ws.Range("J136").Value = 1899.8334
$ws.Range("K136").Value = 107146596
$ws.Range("L136").Value = 5699.5002
$ws.Range("M136").Value = -107144046
$ws.Range("N136").Value = -10799.5002
